# Update example spreadsheet for control monotributistas
#
# The sheet previously had a single "Descarga_MC" style download column
# header left over in K ("Descarga_RCEL"). This change splits the MC
# download status into two columns - "emitidos" (issued) and "recibidos"
# (received) - inserted right before the existing "Descarga_RCEL" column,
# which shifts two columns to the right (K,L -> M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at K/L; the old K (Descarga_RCEL) and its
# data move to M, carrying their formatting with them.
$ws.Range("K1:L1").EntireColumn.Insert()

# New header row values for the inserted columns.
$ws.Range("K1").Value = "Descarga_MC_emitidos"
$ws.Range("L1").Value = "Descarga_MC_recibidos"

# Row 2 (Clave 1 / ABP): was "no" for the single MC-download column,
# now "si" for emitidos and "no" for recibidos.
$ws.Range("K2").Value = "si"
$ws.Range("L2").Value = "no"

# Row 3 (Clave 2 / MB): was "si" for the single MC-download column,
# now "si" for both emitidos and recibidos.
$ws.Range("K3").Value = "si"
$ws.Range("L3").Value = "si"

# Match the author's final selection.
$ws.Application.Goto($ws.Range("L3"))
